$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 7 — the "Updating passwords" security-issue entry.
$ws.Range("B7").Value = "Updating passwords"
$ws.Range("C7").Value = "Login.java, Frame.java, ResetRequest.java, ResetPassword.java"
$ws.Range("D7").Value = "No chance to reset password. If resetting password is allowed and no controls, attackers could reset password whenever wanted"
$ws.Range("E7").Value = "When a password reset is requested, the user must input their email. A reset code (dummy) is sent to their email (dummy). They must then input that reset code, their new password, and their confirmed new password to change their password. "

# The added paragraph text now needs the row tall enough to show it (matches
# the wrapped-text height pattern already used by the other data rows).
$ws.Rows.Item(7).RowHeight = 63

# Leave the selection on the row that was just edited.
$ws.Range("A7").Select()
